$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Remove Incomplete Records")

$data = @{
    7  = @(77.08, 79.69, 71.35, 79.17, 73.96, 75, 72.92, 75.52, 77.6, 76.56)
    8  = @(74.03, 71.86, 68.83, 71.43, 75.32, 76.62, 74.03, 74.46, 76.62, 70.13)
    9  = @(75.09, 77.32, 69.89, 72.86, 78.81, 72.49, 75.46, 68.4, 75.09, 76.58)
    10 = @(74.03, 76.3, 72.73, 75.32, 72.4, 73.38, 72.4, 74.35, 71.43, 73.7)
    11 = @(71.39, 71.39, 76.3, 73.7, 76.3, 74.57, 73.12, 73.99, 71.1, 74.57)
    12 = @(74.48, 68.75, 73.44, 67.19, 74.74, 72.92, 74.22, 74.74, 72.4, 71.61)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 3 + $i   # column C = 3
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

$ws.Range("C13").Select()
